{"js": "// Replace the filename \"PPAML_Challenge_Problem_4.pdf\" (underscore-separated)\n// with \"PPAML-Challenge-Problem-4.pdf\" (hyphen-separated), written as a\n// sequence of separate runs (mirrors the author's incremental retyping),\n// and leave the \"_GoBack\" bookmark wrapped around the new text instead of\n// sitting empty right after it.\n\nconst body = context.document.body;\n\n// Locate the single occurrence of the old filename text.\nconst results = body.search(\"PPAML_Challenge_Problem_4.pdf\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'PPAML_Challenge_Problem_4.pdf' in the document body.\");\n}\n\nconst target = results.items[0];\n\n// Build the replacement as a run-per-chunk fragment (Flat OPC) so each\n// chunk becomes its own <w:r> instead of being coalesced into a single run.\nconst chunks = [\"PPAML-\", \"Challenge\", \"-\", \"Problem\", \"-\", \"4.pdf\"];\nconst runsXml = chunks\n  .map((t) => {\n    // Only mark xml:space=\"preserve\" when the chunk actually needs it\n    // (leading/trailing whitespace) to mirror how Word emits runs.\n    const needsPreserve = /^\\s|\\s$/.test(t);\n    const attr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n    return `<w:r><w:t${attr}>${t}</w:t></w:r>`;\n  })\n  .join(\"\");\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n  runsXml +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\n// Replacing the found range with this fragment drops in six sibling runs at\n// that exact spot, right where the bookmark markers already sit.\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the filename \"PPAML_Challenge_Problem_4.pdf\" (underscore-separated)\n# with \"PPAML-Challenge-Problem-4.pdf\" (hyphen-separated), written as a\n# sequence of separate runs (mirrors the author's incremental retyping),\n# leaving the \"_GoBack\" bookmark wrapped around the new text instead of\n# sitting empty right after it.\n\n$d = $word.ActiveDocument\n\n# Locate the single occurrence of the old filename text.\n$find = $d.Content.Find\n$find.Text = \"PPAML_Challenge_Problem_4.pdf\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find 'PPAML_Challenge_Problem_4.pdf' in the document.\"\n}\n\n# Re-seat a plain Range over the hit's character bounds so InsertXML replaces\n# exactly that span (the Find-owning Range object itself doesn't apply here).\n$target = $d.Range($find.Parent.Start, $find.Parent.End)\n\n# Build the replacement as a run-per-chunk fragment (Flat OPC) so each chunk\n# becomes its own <w:r> instead of being coalesced into a single run.\n$chunks = @(\"PPAML-\", \"Challenge\", \"-\", \"Problem\", \"-\", \"4.pdf\")\n$runsXml = \"\"\nforeach ($chunk in $chunks) {\n    $runsXml += \"<w:r><w:t>$chunk</w:t></w:r>\"\n}\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n       '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n       '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n       '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n       $runsXml +\n       '</w:p></w:body></w:document>' +\n       '</pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($xml)\n"}
